# Append 12 new "foil" rows (SetD/073a.jpg .. SetD/084a.jpg) to the
# retrieval4 worksheet, right after the existing last row (85).
#
# The underlying OOXML diff shows that 12 new shared strings
# ("SetD/073a.jpg" .. "SetD/084a.jpg") were inserted into the shared
# string table right before "targ"/"lure"/"foil", which simply shifts
# those three strings' indices up by 12 (87->99, 88->100, 89->101).
# The actual text values used by the worksheet are unchanged, so we
# just need to add the 12 new rows with the right text values; Excel's
# COM layer takes care of the shared-string bookkeeping itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 86
$lastImageNum = 72

for ($i = 1; $i -le 12; $i++) {
    $row = $startRow + $i - 1
    $imgNum = $lastImageNum + $i
    $imgName = "SetD/{0:D3}a.jpg" -f $imgNum

    $ws.Cells.Item($row, 1).Value = $imgName
    $ws.Cells.Item($row, 2).Value = "foil"
    $ws.Cells.Item($row, 3).Value = 3
}
